# Apply updated team-specific time transition matrix values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value2 = 0.2203389830508475
$ws.Cells.Item(2, 3).Value2 = 0.5169491525423728
$ws.Cells.Item(2, 10).Value2 = 0.008474576271186441
$ws.Cells.Item(2, 16).Value2 = 0.1525423728813559
$ws.Cells.Item(2, 19).Value2 = 0.1016949152542373
# Row 3
$ws.Cells.Item(3, 2).Value2 = 0.01041666666666667
$ws.Cells.Item(3, 3).Value2 = 0.03125
$ws.Cells.Item(3, 10).Value2 = 0.02604166666666667
$ws.Cells.Item(3, 16).Value2 = 0.765625
$ws.Cells.Item(3, 19).Value2 = 0.1666666666666667
# Row 4
$ws.Cells.Item(4, 10).Value2 = 0.01694915254237288
$ws.Cells.Item(4, 16).Value2 = 0.7288135593220338
$ws.Cells.Item(4, 19).Value2 = 0.2542372881355932
# Row 6
$ws.Cells.Item(6, 2).Value2 = 0.04405286343612335
$ws.Cells.Item(6, 4).Value2 = 0.02202643171806168
$ws.Cells.Item(6, 6).Value2 = 0.06607929515418502
$ws.Cells.Item(6, 10).Value2 = 0.3392070484581498
$ws.Cells.Item(6, 15).Value2 = 0.013215859030837
$ws.Cells.Item(6, 17).Value2 = 0.1541850220264317
$ws.Cells.Item(6, 18).Value2 = 0.03524229074889868
$ws.Cells.Item(6, 19).Value2 = 0.3259911894273128
# Row 7
$ws.Cells.Item(7, 2).Value2 = 0.09364548494983277
$ws.Cells.Item(7, 4).Value2 = 0.02675585284280936
$ws.Cells.Item(7, 5).Value2 = 0.003344481605351171
$ws.Cells.Item(7, 6).Value2 = 0.03678929765886288
$ws.Cells.Item(7, 10).Value2 = 0.1237458193979933
$ws.Cells.Item(7, 15).Value2 = 0.04013377926421405
$ws.Cells.Item(7, 17).Value2 = 0.2006688963210702
$ws.Cells.Item(7, 18).Value2 = 0.07023411371237458
$ws.Cells.Item(7, 19).Value2 = 0.4046822742474916
# Row 8
$ws.Cells.Item(8, 2).Value2 = 0.09426229508196721
$ws.Cells.Item(8, 4).Value2 = 0.02254098360655738
$ws.Cells.Item(8, 5).Value2 = 0.002049180327868853
$ws.Cells.Item(8, 6).Value2 = 0.05122950819672131
$ws.Cells.Item(8, 10).Value2 = 0.1147540983606557
$ws.Cells.Item(8, 15).Value2 = 0.02254098360655738
$ws.Cells.Item(8, 17).Value2 = 0.1967213114754098
$ws.Cells.Item(8, 18).Value2 = 0.0778688524590164
$ws.Cells.Item(8, 19).Value2 = 0.4180327868852459
# Row 9
$ws.Cells.Item(9, 2).Value2 = 0.1239316239316239
$ws.Cells.Item(9, 4).Value2 = 0.02991452991452992
$ws.Cells.Item(9, 5).Value2 = 0.004273504273504274
$ws.Cells.Item(9, 6).Value2 = 0.04700854700854701
$ws.Cells.Item(9, 10).Value2 = 0.1282051282051282
$ws.Cells.Item(9, 15).Value2 = 0.02136752136752137
$ws.Cells.Item(9, 17).Value2 = 0.2051282051282051
$ws.Cells.Item(9, 18).Value2 = 0.06837606837606838
$ws.Cells.Item(9, 19).Value2 = 0.3717948717948718
# Row 10
$ws.Cells.Item(10, 2).Value2 = 0.1336134453781513
$ws.Cells.Item(10, 4).Value2 = 0.02605042016806723
$ws.Cells.Item(10, 5).Value2 = 0.001680672268907563
$ws.Cells.Item(10, 6).Value2 = 0.06134453781512605
$ws.Cells.Item(10, 10).Value2 = 0.08991596638655462
$ws.Cells.Item(10, 15).Value2 = 0.02100840336134454
$ws.Cells.Item(10, 17).Value2 = 0.2210084033613445
$ws.Cells.Item(10, 18).Value2 = 0.073109243697479
$ws.Cells.Item(10, 19).Value2 = 0.3722689075630252
# Row 11
$ws.Cells.Item(11, 7).Value2 = 0.1274038461538461
$ws.Cells.Item(11, 10).Value2 = 0.05288461538461538
$ws.Cells.Item(11, 11).Value2 = 0.1634615384615385
$ws.Cells.Item(11, 12).Value2 = 0.6418269230769231
$ws.Cells.Item(11, 19).Value2 = 0.01442307692307692
# Row 12
$ws.Cells.Item(12, 6).Value2 = 0.003448275862068965
$ws.Cells.Item(12, 7).Value2 = 0.7586206896551724
$ws.Cells.Item(12, 10).Value2 = 0.1482758620689655
$ws.Cells.Item(12, 11).Value2 = 0.006896551724137931
$ws.Cells.Item(12, 12).Value2 = 0.05862068965517241
$ws.Cells.Item(12, 19).Value2 = 0.02413793103448276
# Row 13
$ws.Cells.Item(13, 7).Value2 = 0.631578947368421
$ws.Cells.Item(13, 10).Value2 = 0.3157894736842105
$ws.Cells.Item(13, 19).Value2 = 0.05263157894736842
# Row 15
$ws.Cells.Item(15, 6).Value2 = 0.03891050583657588
$ws.Cells.Item(15, 8).Value2 = 0.1439688715953307
$ws.Cells.Item(15, 9).Value2 = 0.0933852140077821
$ws.Cells.Item(15, 10).Value2 = 0.3073929961089494
$ws.Cells.Item(15, 11).Value2 = 0.04280155642023346
$ws.Cells.Item(15, 13).Value2 = 0.01556420233463035
$ws.Cells.Item(15, 15).Value2 = 0.05447470817120623
$ws.Cells.Item(15, 19).Value2 = 0.3035019455252918
# Row 16
$ws.Cells.Item(16, 6).Value2 = 0.03896103896103896
$ws.Cells.Item(16, 8).Value2 = 0.1471861471861472
$ws.Cells.Item(16, 9).Value2 = 0.07792207792207792
$ws.Cells.Item(16, 10).Value2 = 0.3203463203463203
$ws.Cells.Item(16, 11).Value2 = 0.1515151515151515
$ws.Cells.Item(16, 13).Value2 = 0.02164502164502164
$ws.Cells.Item(16, 15).Value2 = 0.03896103896103896
$ws.Cells.Item(16, 19).Value2 = 0.2034632034632035
# Row 17
$ws.Cells.Item(17, 6).Value2 = 0.02173913043478261
$ws.Cells.Item(17, 8).Value2 = 0.199604743083004
$ws.Cells.Item(17, 9).Value2 = 0.116600790513834
$ws.Cells.Item(17, 10).Value2 = 0.3221343873517787
$ws.Cells.Item(17, 11).Value2 = 0.116600790513834
$ws.Cells.Item(17, 13).Value2 = 0.02766798418972332
$ws.Cells.Item(17, 14).Value2 = 0.001976284584980237
$ws.Cells.Item(17, 15).Value2 = 0.06719367588932806
$ws.Cells.Item(17, 19).Value2 = 0.1264822134387352
# Row 18
$ws.Cells.Item(18, 6).Value2 = 0.04191616766467066
$ws.Cells.Item(18, 8).Value2 = 0.2514970059880239
$ws.Cells.Item(18, 9).Value2 = 0.07784431137724551
$ws.Cells.Item(18, 10).Value2 = 0.2994011976047904
$ws.Cells.Item(18, 11).Value2 = 0.1437125748502994
$ws.Cells.Item(18, 13).Value2 = 0.02395209580838323
$ws.Cells.Item(18, 15).Value2 = 0.0718562874251497
$ws.Cells.Item(18, 19).Value2 = 0.08982035928143713
# Row 19
$ws.Cells.Item(19, 6).Value2 = 0.01834862385321101
$ws.Cells.Item(19, 8).Value2 = 0.1961891319689485
$ws.Cells.Item(19, 9).Value2 = 0.08962597035991532
$ws.Cells.Item(19, 10).Value2 = 0.3091037402964009
$ws.Cells.Item(19, 11).Value2 = 0.1496118560338744
$ws.Cells.Item(19, 13).Value2 = 0.02117148906139732
$ws.Cells.Item(19, 14).Value2 = 0.001411432604093155
$ws.Cells.Item(19, 15).Value2 = 0.07198306280875089
$ws.Cells.Item(19, 19).Value2 = 0.1425546930134086
